$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 358
$ws1.Range("F4").Value = 2955
$ws1.Range("F5").Value = 73
$ws1.Range("F6").Value = 618

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F5").Value = 358
$ws4.Range("F6").Value = 2955
$ws4.Range("F7").Value = 73
$ws4.Range("F8").Value = 618
